$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated execution-time results for SP-D (results moved to new directory)
$ws.Range("C2").Value = 114411558.4127764
$ws.Range("D2").Value = 1.828924479393044
$ws.Range("E2").Value = 4437.522373534417
$ws.Range("F2").Value = 4437.522373534417

$ws.Range("C3").Value = 117787547.8801956
$ws.Range("D3").Value = 1.685920957457729
$ws.Range("E3").Value = 4675.947643638125
$ws.Range("F3").Value = 9113.470017172542

$ws.Range("C4").Value = 117794031.8995098
$ws.Range("D4").Value = 1.52173442213647
$ws.Range("E4").Value = 4975.703568426338
$ws.Range("F4").Value = 14089.17358559888

$ws.Range("C5").Value = 117884199.2298288
$ws.Range("D5").Value = 1.631068728650424
$ws.Range("E5").Value = 4638.619952659808
$ws.Range("F5").Value = 18727.79353825869

$ws.Range("C6").Value = 117830475.4317073
$ws.Range("D6").Value = 1.655318423829817
$ws.Range("E6").Value = 5293.955960996308
$ws.Range("F6").Value = 24021.74949925499

$ws.Range("C7").Value = 117716727.8386308
$ws.Range("D7").Value = 1.503755554560929
$ws.Range("E7").Value = 5833.162523811434
$ws.Range("F7").Value = 29854.91202306643

$ws.Range("C8").Value = 116962833.3325183
$ws.Range("D8").Value = 1.664747774119225
$ws.Range("E8").Value = 4221.214802572466
$ws.Range("F8").Value = 34076.12682563889

$ws.Range("C9").Value = 117840207.3609756
$ws.Range("D9").Value = 1.664587833628723
$ws.Range("E9").Value = 3986.156372545263
$ws.Range("F9").Value = 38062.28319818416

$ws.Range("C10").Value = 118333141.7188264
$ws.Range("D10").Value = 1.443363636114659
$ws.Range("E10").Value = 4656.928927635208
$ws.Range("F10").Value = 42719.21212581937
